$wb = $excel.ActiveWorkbook

# Sheets: 1 = test_suite, 2 = AddCustomerTest, 3 = OpenAccountTest
$wsTestSuite = $wb.Worksheets.Item(1)
$wsAddCustomer = $wb.Worksheets.Item(2)

# Data edit: AddCustomerTest!E3 (Sam Tar's runmode) changes from "N" to "Y".
# Once "N" is no longer referenced anywhere it drops out of the shared-strings table.
$wsAddCustomer.Range("E3").Value = "Y"

# Selection/view bookkeeping:
# - AddCustomerTest loses tabSelected and its remembered selection moves to E9.
[void]$wsAddCustomer.Range("E9").Select()

# - test_suite becomes the active/selected tab, with its remembered selection at D9.
#   Selecting the sheet first makes it active, then selecting the cell updates the
#   remembered selection for that sheet (also clearing the workbook's activeTab override).
[void]$wsTestSuite.Select()
[void]$wsTestSuite.Range("D9").Select()
